$d = $word.ActiveDocument

# The document's two flextable tables each have a header row whose top
# and bottom cell borders are drawn in the darker "666666" rule, a second
# row whose top border continues that rule, and a final row whose bottom
# border closes it out. All of those single/666666 cell borders were
# drawn too thick (sz=16, i.e. 2pt) and need to come down to sz=12 (1.5pt).
# Word's Border.LineWidth is expressed in the same "eighth-of-a-point"
# units halved (LineWidth * 2 = sz), so sz=12 == LineWidth 6.
$targetLineWidth = 6

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables($ti)
    $rowCount = $tbl.Rows.Count

    for ($ri = 1; $ri -le $rowCount; $ri++) {
        $row = $tbl.Rows($ri)

        $touchTop = ($ri -eq 1) -or ($ri -eq 2)
        $touchBottom = ($ri -eq 1) -or ($ri -eq $rowCount)

        if (-not ($touchTop -or $touchBottom)) {
            continue
        }

        foreach ($cell in $row.Cells) {
            $borders = $cell.Borders

            if ($touchTop) {
                $top = $borders.Item(-1)
                $top.LineWidth = $targetLineWidth
            }

            if ($touchBottom) {
                $bottom = $borders.Item(-3)
                $bottom.LineWidth = $targetLineWidth
            }
        }
    }
}
